$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 307
$ws.Range("I12").Value = 279.25
$ws.Range("J12").Value = 362.5
$ws.Range("K12").Value = 279.25
$ws.Range("L12").Value = 362.5
$ws.Range("M12").Value = -109.25
$ws.Range("N12").Value = -702.5
$ws.Range("H17").Value = 2814591.5
$ws.Range("J17").Value = 2864849
$ws.Range("L17").Value = 8594547
$ws.Range("N17").Value = -8594883
$ws.Range("H114").Value = 40000
$ws.Range("I114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("M114").ClearContents()
$ws.Range("H116").Value = 7247.1055
$ws.Range("I116").Value = 9515
$ws.Range("K116").Value = 9515
$ws.Range("M116").Value = -6073
$ws.Range("H132").Value = 1786.0526
$ws.Range("I132").Value = 1786.0526
$ws.Range("K132").Value = 5358.1578
$ws.Range("M132").Value = -2828.1578
$ws.Range("H137").Value = 1788.4736
$ws.Range("I137").Value = 1735.7858
$ws.Range("J137").Value = 1936
$ws.Range("K137").Value = 5207.357400000001
$ws.Range("L137").Value = 5808
$ws.Range("M137").Value = -2657.357400000001
$ws.Range("N137").Value = -10908
$ws.Range("H140").Value = 48750
$ws.Range("J140").Value = 48750
$ws.Range("L140").Value = 48750
$ws.Range("N140").Value = -59110
$ws.Range("H141").Value = 4660.5557
$ws.Range("I141").Value = 4660.5557
$ws.Range("K141").Value = 13981.6671
$ws.Range("M141").Value = -8801.667099999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1449.0333
$ws.Range("I2").Value = 1586.2307
$ws.Range("J2").Value = 557.25
$ws.Range("K2").Value = 1586.2307
$ws.Range("L2").Value = 557.25
$ws.Range("M2").Value = -1473.2307
$ws.Range("N2").Value = -783.25
$ws.Range("H32").Value = 5511.982
$ws.Range("I32").Value = 4164.7144
$ws.Range("J32").Value = 9864.691999999999
$ws.Range("K32").Value = 4164.7144
$ws.Range("L32").Value = 9864.691999999999
$ws.Range("M32").Value = -3877.7144
$ws.Range("N32").Value = -10438.692
$ws.Range("H45").Value = 21122.4
$ws.Range("I45").Value = 26028
$ws.Range("J45").Value = 1500
$ws.Range("K45").Value = 26028
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = -25651
$ws.Range("N45").Value = -2254
$ws.Range("H63").Value = 200001940
$ws.Range("I63").Value = 200001940
$ws.Range("K63").Value = 200001940
$ws.Range("M63").Value = -200001254
$ws.Range("H66").Value = 200001940
$ws.Range("I66").Value = 200001940
$ws.Range("K66").Value = 1000009700
$ws.Range("M66").Value = -1000006268
$ws.Range("H74").Value = 1515.6
$ws.Range("I74").Value = 1482.258
$ws.Range("J74").Value = 2204.6667
$ws.Range("K74").Value = 1482.258
$ws.Range("L74").Value = 2204.6667
$ws.Range("M74").Value = -608.258
$ws.Range("N74").Value = -3952.6667
$ws.Range("H77").Value = 1515.6
$ws.Range("I77").Value = 1482.258
$ws.Range("J77").Value = 2204.6667
$ws.Range("K77").Value = 7411.29
$ws.Range("L77").Value = 11023.3335
$ws.Range("M77").Value = -3043.29
$ws.Range("N77").Value = -19759.3335
$ws.Range("H88").Value = 2671.3333
$ws.Range("I88").Value = 3000
$ws.Range("J88").Value = 2507
$ws.Range("K88").Value = 3000
$ws.Range("L88").Value = 2507
$ws.Range("M88").Value = -2594
$ws.Range("N88").Value = -3319
$ws.Range("H91").Value = 2671.3333
$ws.Range("I91").Value = 3000
$ws.Range("J91").Value = 2507
$ws.Range("K91").Value = 3000
$ws.Range("L91").Value = 2507
$ws.Range("M91").Value = -1596
$ws.Range("N91").Value = -5315
$ws.Range("H116").Value = 1449.0333
$ws.Range("I116").Value = 1586.2307
$ws.Range("J116").Value = 557.25
$ws.Range("K116").Value = 1586.2307
$ws.Range("L116").Value = 557.25
$ws.Range("M116").Value = 707.7692999999999
$ws.Range("N116").Value = -5145.25
$ws.Range("H132").Value = 3709
$ws.Range("I132").Value = 1327.7646
$ws.Range("K132").Value = 3983.2938
$ws.Range("M132").Value = -1453.2938
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1449.0333
$ws.Range("I3").Value = 1586.2307
$ws.Range("J3").Value = 557.25
$ws.Range("K3").Value = 1586.2307
$ws.Range("L3").Value = 557.25
$ws.Range("M3").Value = -1472.2307
$ws.Range("N3").Value = -785.25
$ws.Range("H99").Value = 83334420
$ws.Range("I99").Value = 90910130
$ws.Range("K99").Value = 90910130
$ws.Range("M99").Value = -90908632
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H140").Value = 45972.637
$ws.Range("J140").Value = 45972.637
$ws.Range("L140").Value = 45972.637
$ws.Range("N140").Value = -56332.637
$ws.Range("H141").Value = 40780
$ws.Range("J141").Value = 40780
$ws.Range("L141").Value = 40780
$ws.Range("N141").Value = -51140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1382
$ws.Range("I19").Value = 303.33334
$ws.Range("J19").Value = 3000
$ws.Range("K19").Value = 303.33334
$ws.Range("L19").Value = 3000
$ws.Range("M19").Value = -133.33334
$ws.Range("N19").Value = -3340
$ws.Range("H24").Value = 1382
$ws.Range("I24").Value = 303.33334
$ws.Range("J24").Value = 3000
$ws.Range("K24").Value = 303.33334
$ws.Range("L24").Value = 3000
$ws.Range("M24").Value = -133.33334
$ws.Range("N24").Value = -3340
$ws.Range("H31").Value = 9098.571
$ws.Range("I31").Value = 2356.7144
$ws.Range("J31").Value = 22582.285
$ws.Range("K31").Value = 2356.7144
$ws.Range("L31").Value = 22582.285
$ws.Range("M31").Value = -2061.7144
$ws.Range("N31").Value = -23172.285
$ws.Range("H34").Value = 9098.571
$ws.Range("I34").Value = 2356.7144
$ws.Range("J34").Value = 22582.285
$ws.Range("K34").Value = 2356.7144
$ws.Range("L34").Value = 22582.285
$ws.Range("M34").Value = -2154.7144
$ws.Range("N34").Value = -22986.285
$ws.Range("H58").Value = 1545
$ws.Range("I58").Value = 1112.25
$ws.Range("J58").Value = 1915.9286
$ws.Range("K58").Value = 1112.25
$ws.Range("L58").Value = 1915.9286
$ws.Range("M58").Value = -909.25
$ws.Range("N58").Value = -2321.9286
$ws.Range("H134").Value = 3196.926
$ws.Range("I134").Value = 4040.5881
$ws.Range("J134").Value = 1762.7
$ws.Range("K134").Value = 12121.7643
$ws.Range("L134").Value = 5288.1
$ws.Range("M134").Value = -9586.764299999999
$ws.Range("N134").Value = -10358.1
$ws.Range("H136").Value = 1545
$ws.Range("I136").Value = 1112.25
$ws.Range("J136").Value = 1915.9286
$ws.Range("K136").Value = 3336.75
$ws.Range("L136").Value = 5747.7858
$ws.Range("M136").Value = -786.75
$ws.Range("N136").Value = -10847.7858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 4000626
$ws.Range("I113").Value = 12500412
$ws.Range("J113").Value = 909794.5600000001
$ws.Range("K113").Value = 37501236
$ws.Range("L113").Value = 2729383.68
$ws.Range("M113").Value = -37499066
$ws.Range("N113").Value = -2733723.68
$ws.Range("H129").Value = 2267.3845
$ws.Range("I129").Value = 1582.7142
$ws.Range("J129").Value = 3066.1667
$ws.Range("K129").Value = 4748.142599999999
$ws.Range("L129").Value = 9198.500100000001
$ws.Range("M129").Value = 251.8574000000008
$ws.Range("N129").Value = -19198.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1229
$ws.Range("I16").Value = 827.2857
$ws.Range("J16").Value = 2166.3333
$ws.Range("K16").Value = 827.2857
$ws.Range("L16").Value = 2166.3333
$ws.Range("M16").Value = -657.2857
$ws.Range("N16").Value = -2506.3333
$ws.Range("H68").Value = 100002320
$ws.Range("I68").Value = 1854.5
$ws.Range("J68").Value = 166669300
$ws.Range("K68").Value = 1854.5
$ws.Range("L68").Value = 166669300
$ws.Range("M68").Value = -1105.5
$ws.Range("N68").Value = -166670798
$ws.Range("H71").Value = 100002320
$ws.Range("I71").Value = 1854.5
$ws.Range("J71").Value = 166669300
$ws.Range("K71").Value = 9272.5
$ws.Range("L71").Value = 833346500
$ws.Range("M71").Value = -5528.5
$ws.Range("N71").Value = -833353988
$ws.Range("H140").Value = 44920.875
$ws.Range("J140").Value = 44920.875
$ws.Range("L140").Value = 44920.875
$ws.Range("N140").Value = -55280.875
$ws.Range("H141").Value = 49657.145
$ws.Range("J141").Value = 49657.145
$ws.Range("L141").Value = 49657.145
$ws.Range("N141").Value = -60017.145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 54400
$ws.Range("J46").Value = 54400
$ws.Range("L46").Value = 54400
$ws.Range("N46").Value = -54862
$ws.Range("H86").Value = 24714.285
$ws.Range("J86").Value = 24714.285
$ws.Range("L86").Value = 24714.285
$ws.Range("N86").Value = -26960.285
$ws.Range("H89").Value = 24714.285
$ws.Range("J89").Value = 24714.285
$ws.Range("L89").Value = 123571.425
$ws.Range("N89").Value = -134803.425
$ws.Range("H100").Value = 7059.467
$ws.Range("I100").Value = 8682.916999999999
$ws.Range("J100").Value = 565.6667
$ws.Range("K100").Value = 17365.834
$ws.Range("L100").Value = 1131.3334
$ws.Range("M100").Value = -16824.834
$ws.Range("N100").Value = -2213.3334
$ws.Range("H134").Value = 54400
$ws.Range("J134").Value = 54400
$ws.Range("L134").Value = 163200
$ws.Range("N134").Value = -168270
